$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that originally sat right after the
#    "TONG CONG" run (it will be re-created further below, in its new spot).
# ---------------------------------------------------------------------------
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
}

# ---------------------------------------------------------------------------
# 2) Find the empty placeholder paragraph (style "List") that sits right
#    after the valuation table and fill it in with the "Gia tri dinh gia"
#    line, followed by a brand new "Bang chu" paragraph.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$anchorIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "*{ID6_TOTAL_GT_CHO_VAY}*") {
        $anchorIndex = $i
        break
    }
}

$emptyIndex = -1
for ($j = $anchorIndex + 1; $j -le $paras.Count; $j++) {
    $candidate = $paras.Item($j)
    if ($candidate.Style.NameLocal -eq "List" -and $candidate.Range.Text.Trim() -eq "") {
        $emptyIndex = $j
        break
    }
}

$emptyPara = $paras.Item($emptyIndex)
$insertionRange = $emptyPara.Range
$insertionRange.Collapse(0)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:pStyle w:val="List"/><w:tabs><w:tab w:val="left" w:pos="567"/></w:tabs><w:spacing w:after="0" w:line="276" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t>Giá</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t>trị</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t>định</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t>giá</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t xml:space="preserve"> : </w:t></w:r><w:bookmarkStart w:id="100" w:name="_GoBack"/><w:r><w:rPr><w:b/></w:rPr><w:t>{ID6_VALUE}</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>đồng</w:t></w:r><w:bookmarkEnd w:id="100"/><w:proofErr w:type="spellEnd"/></w:p>
<w:p><w:pPr><w:pStyle w:val="List"/><w:tabs><w:tab w:val="left" w:pos="567"/></w:tabs><w:spacing w:after="0" w:line="276" w:lineRule="auto"/><w:rPr><w:i w:val="0"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Bằng</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>chữ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> : </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>{ID6_TEXT}</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionRange.InsertXML($xml)
